# Auto-generated edit script applying scheduled market-data refresh values
# to the Alpha_Profits workbook (per-sheet currentAveragePrice / LevePrice / LeveProfit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5349.6113
$ws.Range("I64").Value = 5339.1
$ws.Range("J64").Value = 5362.75
$ws.Range("K64").Value = 5339.1
$ws.Range("L64").Value = 5362.75
$ws.Range("M64").Value = -5091.1
$ws.Range("N64").Value = -5858.75

$ws.Range("H67").Value = 5349.6113
$ws.Range("I67").Value = 5339.1
$ws.Range("J67").Value = 5362.75
$ws.Range("K67").Value = 5339.1
$ws.Range("L67").Value = 5362.75
$ws.Range("M67").Value = -4481.1
$ws.Range("N67").Value = -7078.75

$ws.Range("H103").Value = 611.80646
$ws.Range("I103").Value = 604.04346
$ws.Range("K103").Value = 1812.13038
$ws.Range("M103").Value = -1226.13038

$ws.Range("H106").Value = 2699.7778
$ws.Range("I106").Value = 2699.7778
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2699.7778
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").Value = -2068.7778

$ws.Range("H113").Value = 4467.1875
$ws.Range("I113").Value = 2566.1428
$ws.Range("K113").Value = 2566.1428
$ws.Range("M113").Value = 687.8571999999999

$ws.Range("H137").Value = 2186.0344
$ws.Range("I137").Value = 1622.125
$ws.Range("J137").Value = 2400.8572
$ws.Range("K137").Value = 4866.375
$ws.Range("L137").Value = 7202.571599999999
$ws.Range("M137").Value = -2316.375
$ws.Range("N137").Value = -12302.5716

$ws.Range("H138").Value = 3189.743
$ws.Range("I138").Value = 1316
$ws.Range("J138").Value = 3431.516
$ws.Range("K138").Value = 3948
$ws.Range("L138").Value = 10294.548
$ws.Range("M138").Value = 1192
$ws.Range("N138").Value = -20574.548

$ws.Range("H141").Value = 195398.6
$ws.Range("I141").Value = 322331.34
$ws.Range("J141").Value = 4999.5
$ws.Range("K141").Value = 966994.02
$ws.Range("L141").Value = 14998.5
$ws.Range("M141").Value = -961814.02
$ws.Range("N141").Value = -25358.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 11049.667
$ws.Range("I31").Value = 11049.667
$ws.Range("K31").Value = 11049.667
$ws.Range("M31").Value = -10755.667

$ws.Range("H32").Value = 3647971.8
$ws.Range("I32").Value = 3850892.5
$ws.Range("K32").Value = 3850892.5
$ws.Range("M32").Value = -3850605.5

$ws.Range("H40").Value = 32966.332
$ws.Range("I40").Value = 24450
$ws.Range("J40").Value = 49999
$ws.Range("K40").Value = 24450
$ws.Range("L40").Value = 49999
$ws.Range("M40").Value = -24274
$ws.Range("N40").Value = -50351

$ws.Range("H60").Value = 16949.5
$ws.Range("I60").Value = 16949.5
$ws.Range("K60").Value = 16949.5
$ws.Range("M60").Value = -16216.5

$ws.Range("H61").Value = 2408.375
$ws.Range("I61").Value = 1377
$ws.Range("J61").Value = 5502.5
$ws.Range("K61").Value = 1377
$ws.Range("L61").Value = 5502.5
$ws.Range("M61").Value = -1165
$ws.Range("N61").Value = -5926.5

$ws.Range("H74").Value = 2329.8667
$ws.Range("I74").Value = 1575.3572
$ws.Range("K74").Value = 1575.3572
$ws.Range("M74").Value = -701.3571999999999

$ws.Range("H77").Value = 2329.8667
$ws.Range("I77").Value = 1575.3572
$ws.Range("K77").Value = 7876.786
$ws.Range("M77").Value = -3508.786

$ws.Range("H132").Value = 2688.0454
$ws.Range("I132").Value = 2016.8823
$ws.Range("J132").Value = 4970
$ws.Range("K132").Value = 6050.6469
$ws.Range("L132").Value = 14910
$ws.Range("M132").Value = -3520.6469
$ws.Range("N132").Value = -19970

$ws.Range("H136").Value = 2408.375
$ws.Range("I136").Value = 1377
$ws.Range("J136").Value = 5502.5
$ws.Range("K136").Value = 4131
$ws.Range("L136").Value = 16507.5
$ws.Range("M136").Value = -1581
$ws.Range("N136").Value = -21607.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3594.738
$ws.Range("I86").Value = 2255.158
$ws.Range("J86").Value = 4701.3477
$ws.Range("K86").Value = 2255.158
$ws.Range("L86").Value = 4701.3477
$ws.Range("M86").Value = -1132.158
$ws.Range("N86").Value = -6947.3477

$ws.Range("H89").Value = 3594.738
$ws.Range("I89").Value = 2255.158
$ws.Range("J89").Value = 4701.3477
$ws.Range("K89").Value = 11275.79
$ws.Range("L89").Value = 23506.7385
$ws.Range("M89").Value = -5659.789999999999
$ws.Range("N89").Value = -34738.7385

$ws.Range("H134").Value = 41671844
$ws.Range("I134").Value = 5913.857
$ws.Range("J134").Value = 333333340
$ws.Range("K134").Value = 17741.571
$ws.Range("L134").Value = 1000000020
$ws.Range("M134").Value = -15206.571
$ws.Range("N134").Value = -1000005090

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2501.9375
$ws.Range("I99").Value = 2043.1666
$ws.Range("K99").Value = 2043.1666
$ws.Range("M99").Value = -545.1666

$ws.Range("H107").Value = 1543.96
$ws.Range("J107").Value = 2651.2727
$ws.Range("L107").Value = 2651.2727
$ws.Range("N107").Value = -6491.2727

$ws.Range("H122").Value = 4177.048
$ws.Range("I122").Value = 4326.684
$ws.Range("J122").Value = 2755.5
$ws.Range("K122").Value = 12980.052
$ws.Range("L122").Value = 8266.5
$ws.Range("M122").Value = -10530.052
$ws.Range("N122").Value = -13166.5

$ws.Range("H126").Value = 2501.9375
$ws.Range("I126").Value = 2043.1666
$ws.Range("K126").Value = 6129.4998
$ws.Range("M126").Value = -3659.4998

$ws.Range("H134").Value = 5002709.5
$ws.Range("I134").Value = 1933.4667
$ws.Range("J134").Value = 20005038
$ws.Range("K134").Value = 5800.4001
$ws.Range("L134").Value = 60015114
$ws.Range("M134").Value = -3265.4001
$ws.Range("N134").Value = -60020184

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 10101010
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H27").Value = 10101010
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H55").Value = 4319.4
$ws.Range("I55").Value = 2931
$ws.Range("K55").Value = 8793
$ws.Range("M55").Value = -8616

$ws.Range("H94").Value = 13997
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").ClearContents()
$ws.Range("N96").Value = 0

$ws.Range("H98").Value = 2773.75
$ws.Range("I98").Value = 2697.5
$ws.Range("J98").Value = 2850
$ws.Range("K98").Value = 8092.5
$ws.Range("L98").Value = 8550
$ws.Range("M98").Value = -6594.5
$ws.Range("N98").Value = -11546

$ws.Range("H99").Value = 1500
$ws.Range("I99").Value = 1500
$ws.Range("K99").Value = 4500
$ws.Range("M99").Value = -2254

$ws.Range("H103").Value = 514.1667
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4758

$ws.Range("H104").Value = 3432.6365
$ws.Range("I104").Value = 3417.6667
$ws.Range("J104").Value = 3500
$ws.Range("K104").Value = 10253.0001
$ws.Range("L104").Value = 10500
$ws.Range("M104").Value = -7632.000100000001
$ws.Range("N104").Value = -15742

$ws.Range("H107").Value = 837.2174
$ws.Range("I107").Value = 325.33334
$ws.Range("J107").Value = 1166.2858
$ws.Range("K107").Value = 976.0000200000001
$ws.Range("L107").Value = 3498.8574
$ws.Range("M107").Value = 943.9999799999999
$ws.Range("N107").Value = -7338.857400000001

$ws.Range("H109").Value = 1000
$ws.Range("I109").Value = 1000
$ws.Range("J109").Value = 0
$ws.Range("K109").Value = 3000
$ws.Range("L109").ClearContents()
$ws.Range("M109").Value = -1960
$ws.Range("N109").Value = 0

$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").ClearContents()

$ws.Range("H112").Value = 3527
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H113").Value = 484.125
$ws.Range("J113").Value = 625.625
$ws.Range("L113").Value = 1876.875
$ws.Range("N113").Value = -6216.875

$ws.Range("H137").Value = 6270.65
$ws.Range("I137").Value = 3406.875
$ws.Range("K137").Value = 10220.625
$ws.Range("M137").Value = -5120.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1680.1177
$ws.Range("I132").Value = 1599
$ws.Range("K132").Value = 4797
$ws.Range("M132").Value = -2267

$ws.Range("H136").Value = 53370.8
$ws.Range("J136").Value = 53370.8
$ws.Range("L136").Value = 160112.4
$ws.Range("N136").Value = -165212.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1560.7916
$ws.Range("I16").Value = 1322
$ws.Range("J16").Value = 3232.3333
$ws.Range("K16").Value = 1322
$ws.Range("L16").Value = 3232.3333
$ws.Range("M16").Value = -1152
$ws.Range("N16").Value = -3572.3333

$ws.Range("H68").Value = 4099.5
$ws.Range("I68").Value = 3749.375
$ws.Range("J68").Value = 5500
$ws.Range("K68").Value = 3749.375
$ws.Range("L68").Value = 5500
$ws.Range("M68").Value = -3000.375
$ws.Range("N68").Value = -6998

$ws.Range("H71").Value = 4099.5
$ws.Range("I71").Value = 3749.375
$ws.Range("J71").Value = 5500
$ws.Range("K71").Value = 18746.875
$ws.Range("L71").Value = 27500
$ws.Range("M71").Value = -15002.875
$ws.Range("N71").Value = -34988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4800
$ws.Range("I62").Value = 4863.636
$ws.Range("J62").Value = 4625
$ws.Range("K62").Value = 4863.636
$ws.Range("L62").Value = 4625
$ws.Range("M62").Value = -4239.636
$ws.Range("N62").Value = -5873

$ws.Range("H65").Value = 4800
$ws.Range("I65").Value = 4863.636
$ws.Range("J65").Value = 4625
$ws.Range("K65").Value = 24318.18
$ws.Range("L65").Value = 23125
$ws.Range("M65").Value = -21198.18
$ws.Range("N65").Value = -29365

$ws.Range("H132").Value = 3997.8572
$ws.Range("I132").Value = 3920.7693
$ws.Range("K132").Value = 11762.3079
$ws.Range("M132").Value = -9232.3079

